# Commit: "Hi Ramesh commiting with the logout action and fetching the URL"
#
# This script:
#  1. Rebuilds the "Logout" sheet so it has the same header row as "Login"
#     plus two test-data rows (reusing Login's row-2 formatting/values),
#     renames the first column values to the new test names, and fetches
#     (hyperlinks) the "nd1432@blr" value in column E for both rows.
#  2. Adds a new blank worksheet named "Sheet4" after "Template".
#  3. Updates view state (zoom 65%, selections) on every sheet and makes
#     "Logout" the active tab.

$wb = $excel.ActiveWorkbook

$login  = $wb.Worksheets.Item("Login")
$logout = $wb.Worksheets.Item("Logout")
$tmpl   = $wb.Worksheets.Item("Template")

# --- 1. Rebuild Logout sheet contents -------------------------------------

# Bring in the header row + first data row layout/styling from Login.
$login.Range("A1:H1").Copy($logout.Range("A1:H1"))
$login.Range("A2:H2").Copy($logout.Range("A2:H2"))
$login.Range("A2:H2").Copy($logout.Range("A3:H3"))

# New test-name descriptions for the two scenario rows.
$logout.Range("A2").Value = "validate login feature with login"
$logout.Range("A3").Value = "Validate with logout"

# Replace the old "Shalini@Navadhiti.Com" hyperlink with fresh hyperlinks
# that fetch/point at "nd1432@blr" for both rows.
foreach ($h in @($logout.Hyperlinks)) {
    $h.Delete()
}
$logout.Hyperlinks.Add($logout.Range("E2"), "mailto:nd1432@blr", "", "", "nd1432@blr")
$logout.Hyperlinks.Add($logout.Range("E3"), "mailto:nd1432@blr", "", "", "nd1432@blr")

# --- 2. Add new "Sheet4" worksheet -----------------------------------------

$tmpl.Copy($null, $tmpl)
$sheet4 = $wb.Worksheets.Item($tmpl.Index() + 1)
$sheet4.Name = "Sheet4"
$sheet4.Rows("1:2").Delete()

# --- 3. View state: zoom 65% everywhere, per-sheet selections --------------

$login.Activate()
$excel.ActiveWindow.Zoom = 65
$login.Range("A1").Select()

$sheet4.Activate()
$excel.ActiveWindow.Zoom = 65
$sheet4.Range("A1").Select()

$tmpl.Activate()
$excel.ActiveWindow.Zoom = 65
$tmpl.Range("P24").Select()

$logout.Activate()
$excel.ActiveWindow.Zoom = 65
$logout.Range("H3").Select()

Write-Host "Logout action + URL fetch changes applied"
